# Add a comment linking to a YouTube video, on slide 1 (the title slide),
# authored by the existing "Tom Lever" comment author (reused from
# ppt/commentAuthors.xml, id="1"). This bumps that author's lastIdx from 1
# to 2, and creates ppt/comments/comment1.xml with a single p:cm (idx="2")
# holding the link text.
#
# PowerPoint's Comments.Add(Left, Top, Author, AuthorInitials, Text) takes
# Left/Top in points and stores them as EMUs (x12700) in the underlying
# p:pos element. The target position in the saved OOXML is x="10" y="10"
# (i.e. 10 EMU, not 10pt), so Left/Top are passed as 10/12700 point so
# that, once multiplied back out to EMU on save, we land exactly on
# x="10" y="10".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$emuToPoint = 10.0 / 12700.0

$comment = $s.Comments.Add($emuToPoint, $emuToPoint, "Tom Lever", "TL", "https://www.youtube.com/watch?v=PjYQt71666g")
